$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Su (t/sq.m.)" column (G) values for rows 4-8
$ws.Range("G4").Value = 23
$ws.Range("G5").Value = 22
$ws.Range("G6").Value = 21
$ws.Range("G7").Value = 20
$ws.Range("G8").Value = 18

# Move the active selection from G4 to G9
$ws.Range("G9").Select()
